$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rename json file paths: insert "unif_filtres/" right after "jsons_train/"
# for both the "before change" (column A) and "after change" (column B) paths.

$ws.Range("A2").Value2 = "jsons_train/unif_filtres/IT__web.json"
$ws.Range("B2").Value2 = "jsons_train/unif_filtres/IT__web_same_as_browser.json"

$ws.Range("A3").Value2 = "jsons_train/unif_filtres/IT__web.json"
$ws.Range("B3").Value2 = "jsons_train/unif_filtres/IT__web_same_as_category.json"

$ws.Range("A4").Value2 = "jsons_train/unif_filtres/IT__web.json"
$ws.Range("B4").Value2 = "jsons_train/unif_filtres/IT__web_same_as_region.json"

$ws.Range("A5").Value2 = "jsons_train/unif_filtres/IT__web_bis.json"
$ws.Range("B5").Value2 = "jsons_train/unif_filtres/IT__web_bis_same_as_country.json"

$ws.Range("A6").Value2 = "jsons_train/unif_filtres/IT__web_bis.json"
$ws.Range("B6").Value2 = "jsons_train/unif_filtres/IT__web_bis_same_as_region.json"

$ws.Range("A7").Value2 = "jsons_train/unif_filtres/IT__web_bis.json"
$ws.Range("B7").Value2 = "jsons_train/unif_filtres/IT__web_bis_same_as_browser.json"

$ws.Range("A8").Value2 = "jsons_train/unif_filtres/IT__web_bis.json"
$ws.Range("B8").Value2 = "jsons_train/unif_filtres/IT__web_bis_browser_same_as_region.json"

$ws.Range("A9").Value2 = "jsons_train/unif_filtres/airbnb.json"
$ws.Range("B9").Value2 = "jsons_train/unif_filtres/airbnb_same_as_cancellation.json"

$ws.Range("A10").Value2 = "jsons_train/unif_filtres/airbnb.json"
$ws.Range("B10").Value2 = "jsons_train/unif_filtres/airbnb_same_as_quartier.json"

$ws.Range("A11").Value2 = "jsons_train/unif_filtres/airbnb.json"
$ws.Range("B11").Value2 = "jsons_train/unif_filtres/airbnb_same_as_room_type.json"

# Update the selected cell in the sheet view (was C7, now C17)
$ws.Range("C17").Select()
